# Add the team's season record (Wins / Losses / Ties) to the sheet.
# These three new columns follow the existing "Unnamed: 28" column (AC)
# and are populated with the same record for every player row, since the
# record belongs to the team/season rather than the individual player.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled like the rest of row 1 (bold, bordered,
# centered horizontally, top-aligned vertically).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous

# Season record for every player row (2 through 56): 86 wins, 76 losses,
# 0 ties.
$firstRow = 2
$lastRow = 56
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 86  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 76  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
